# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (used only by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"     (used by the Slide Master / all slides)
#
# The authored change swaps the two themes' contents (color scheme) so that the
# deck's slide master now carries the "Office Theme" palette while the notes
# master ends up with the former "Integral" palette (fontScheme/fmtScheme are
# identical between the two themes, only clrScheme differs).
#
# The PowerPoint object model only exposes the *slide-facing* theme's color
# slots (via Slide.ThemeColorScheme), so we push the target ("Office Theme")
# RGB values into that scheme, in clrScheme order:
#   dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function ToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Length; $i++) {
    $tcs.Item($i + 1).RGB = ToRGB($officeThemeColors[$i])
}
